# إضافة حدث جديد في Card19 by admin at 2025-12-08 09:29:39
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card19")

# --- Row 18: the previously-empty B18:K18 cells get the literal text "nan",
#     matching the pattern used throughout the rest of the sheet. ---
$ws.Range("B18:K18").Value = "nan"

# --- Row 19: brand new service event row for Card19. ---
# A19 holds the card number as text ("19"), same as the rest of column A.
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = "19"
$ws.Range("A19").Style = "Normal"

# B19:K19 stay blank (present but empty cells), same as B18:K18 were before this edit.
$ws.Range("B19:K19").Style = "Normal"

$ws.Range("L19").Value = "26\10\2025"
$ws.Range("M19").Value = "902 t"
$ws.Range("N19").Value = "تم تغيير الجرائد الخلفيه (1_5_8)"
$ws.Range("O19").Value = "الخبير"
